$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '59.149.22'

$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +0.46%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.611.49'

$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -0.37%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'

$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '560.26'

$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +5.15%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '144.00'

$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +0.40%  '

$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -0.20%  '

$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +5.01%  '

$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -1.36%  '

$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +0.24%  '

$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +5.55%  '

$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +0.10%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.073.56'

$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -0.46%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '59.108.68'

$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +0.52%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '21.15'

$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +1.06%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.625.55'

$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +0.04%  '

$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -0.18%  '

$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +1.40%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '338.07'

$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +0.62%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '10.15'

$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -0.12%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.19'

$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -0.39%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.998'

$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -0.10%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '66.07'

$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -0.50%  '

$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +3.76%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.165'

$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +1.10%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.995'

$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -0.50%  '

$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +0.09%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.0₃0766'

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.998'

$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -0.12%  '

$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +3.22%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.03'

$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +3.52%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '154.81'

$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +2.62%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '18.93'

$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +1.29%  '

$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +1.76%  '

$ws.Range('B35').Value = 'Fetch.AI'

$ws.Range('C35').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.914'

$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +11.05%  '

$ws.Range('B36').Value = 'SuiNetwork'

$ws.Range('C36').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.908'

$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +9.20%  '

$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +1.73%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '37.22'

$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +0.31%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.47'

$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +3.13%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.61'

$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +1.32%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '285.87'

$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +1.49%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.997'

$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -0.26%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.601'

$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +1.24%  '

$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +1.38%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0958'

$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +2.28%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '10.63'

$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -0.93%  '

$ws.Range('B47').Value = 'RenderToken'

$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.69'

$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +3.98%  '

$ws.Range('B48').Value = 'VeChain'

$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0228'

$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +1.71%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.954.97'

$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +0.65%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '118.04'

$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +5.84%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '18.11'

$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +0.53%  '
